$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: Volume number 40 -> 41, and week date range ---
# "Volume 32   Number  40" -> replace the trailing "40" (chars 21-22, 1-based) with "41"
$ws.Range("A8").Characters(21, 2).Text = "41"
# "Report Covering the Week  9/29/2025  Through  10/5/2025"
# replace "9/29/2025" (chars 27-35) with "10/6/2025" (same length, so the 2nd date offset is unaffected)
$ws.Range("C9").Characters(27, 9).Text = "10/6/2025"
# replace "10/5/2025" (chars 47-55) with "10/12/2025"
$ws.Range("C9").Characters(47, 9).Text = "10/12/2025"

# --- Plain numeric cell updates (values recomputed for the new reporting week) ---
$ws.Range("N15").Value2 = 10.526315789473
$ws.Range("C16").Value2 = 2
$ws.Range("D16").Value2 = 2
$ws.Range("E16").Value2 = 0
$ws.Range("F16").Value2 = 15
$ws.Range("H16").Value2 = 36.363636363636
$ws.Range("I16").Value2 = 94
$ws.Range("J16").Value2 = 115
$ws.Range("K16").Value2 = -18.260869565217
$ws.Range("L16").Value2 = -23.577235772357
$ws.Range("M16").Value2 = -24.8
$ws.Range("N16").Value2 = -79.828326180257
$ws.Range("C17").Value2 = 1
$ws.Range("D17").Value2 = 2
$ws.Range("E17").Value2 = -50
$ws.Range("F17").Value2 = 10
$ws.Range("G17").Value2 = 19
$ws.Range("H17").Value2 = -47.368421052631
$ws.Range("I17").Value2 = 179
$ws.Range("J17").Value2 = 194
$ws.Range("K17").Value2 = -7.731958762886
$ws.Range("L17").Value2 = 8.484848484848
$ws.Range("M17").Value2 = 90.425531914893
$ws.Range("N17").Value2 = -20.089285714285
$ws.Range("C18").Value2 = 3
$ws.Range("E18").Value2 = -25
$ws.Range("F18").Value2 = 8
$ws.Range("G18").Value2 = 12
$ws.Range("H18").Value2 = -33.333333333333
$ws.Range("I18").Value2 = 108
$ws.Range("J18").Value2 = 132
$ws.Range("K18").Value2 = -18.181818181818
$ws.Range("L18").Value2 = -6.086956521739
$ws.Range("M18").Value2 = -15.625
$ws.Range("N18").Value2 = -88.039867109634
$ws.Range("C19").Value2 = 10
$ws.Range("D19").Value2 = 15
$ws.Range("E19").Value2 = -33.333333333333
$ws.Range("F19").Value2 = 36
$ws.Range("G19").Value2 = 48
$ws.Range("H19").Value2 = -25
$ws.Range("I19").Value2 = 447
$ws.Range("J19").Value2 = 606
$ws.Range("K19").Value2 = -26.237623762376
$ws.Range("L19").Value2 = -15.500945179584
$ws.Range("M19").Value2 = 67.415730337078
$ws.Range("N19").Value2 = 24.860335195530
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 9
$ws.Range("E20").Value2 = -55.555555555555
$ws.Range("F20").Value2 = 19
$ws.Range("H20").Value2 = -24
$ws.Range("I20").Value2 = 216
$ws.Range("J20").Value2 = 263
$ws.Range("K20").Value2 = -17.870722433460
$ws.Range("L20").Value2 = -20.879120879120
$ws.Range("M20").Value2 = 125
$ws.Range("N20").Value2 = -85.041551246537
$ws.Range("D21").Value2 = 32
$ws.Range("E21").Value2 = -37.5
$ws.Range("F21").Value2 = 90
$ws.Range("G21").Value2 = 118
$ws.Range("H21").Value2 = -23.728813559322
$ws.Range("I21").Value2 = 1066
$ws.Range("J21").Value2 = 1322
$ws.Range("K21").Value2 = -19.364599092284
$ws.Range("L21").Value2 = -12.479474548440
$ws.Range("M21").Value2 = 47.645429362880
$ws.Range("N21").Value2 = -68.839520608009
$ws.Range("L22").Value2 = -12.5
$ws.Range("M22").Value2 = -53.333333333333
$ws.Range("F23").Value2 = 3
$ws.Range("G23").Value2 = 2
$ws.Range("H23").Value2 = 50
$ws.Range("I23").Value2 = 43
$ws.Range("K23").Value2 = -15.686274509803
$ws.Range("L23").Value2 = -33.846153846153
$ws.Range("M23").Value2 = 30.303030303030
$ws.Range("C24").Value2 = 15
$ws.Range("D24").Value2 = 25
$ws.Range("E24").Value2 = -40
$ws.Range("F24").Value2 = 103
$ws.Range("H24").Value2 = 18.390804597701
$ws.Range("I24").Value2 = 891
$ws.Range("J24").Value2 = 868
$ws.Range("K24").Value2 = 2.649769585253
$ws.Range("L24").Value2 = 2.296211251435
$ws.Range("M24").Value2 = 39.436619718309
$ws.Range("C25").Value2 = 4
$ws.Range("D25").Value2 = 8
$ws.Range("E25").Value2 = -50
$ws.Range("F25").Value2 = 30
$ws.Range("H25").Value2 = 36.363636363636
$ws.Range("I25").Value2 = 257
$ws.Range("J25").Value2 = 334
$ws.Range("K25").Value2 = -23.053892215568
$ws.Range("L25").Value2 = -24.853801169590
$ws.Range("C26").Value2 = 5
$ws.Range("D26").Value2 = 6
$ws.Range("E26").Value2 = -16.666666666666
$ws.Range("F26").Value2 = 33
$ws.Range("H26").Value2 = -25
$ws.Range("I26").Value2 = 302
$ws.Range("J26").Value2 = 274
$ws.Range("K26").Value2 = 10.218978102189
$ws.Range("L26").Value2 = 21.774193548387
$ws.Range("M26").Value2 = -1.307189542483
$ws.Range("F28").Value2 = 2
$ws.Range("H28").Value2 = -33.333333333333
$ws.Range("J28").Value2 = 39
$ws.Range("K28").Value2 = -15.384615384615
$ws.Range("L28").Value2 = 32
$ws.Range("L29").Value2 = 80
$ws.Range("L30").Value2 = 20

# --- Cells that flip from a number to the text placeholders "0" / "***.*" ---
# Setting NumberFormat to Text ("@") first forces the digit-looking string "0" to
# actually store as text instead of being reinterpreted as the number 0.
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "***.*"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "***.*"

# Re-apply the same format (right-aligned, General) used by the other text cells in
# these rows (e.g. C15/C27) so the style index matches instead of keeping the ad-hoc "@" style.
$ws.Range("C15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)
$ws.Range("C27").Copy()
$ws.Range("D27:E27").PasteSpecial(-4122)

# --- Cell that flips from the text placeholder "0" back to a real number ---
$ws.Range("F23").Copy()  # F23 already uses the plain integer style for this row
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value2 = 1
